$d = $word.ActiveDocument

$opmaak = $d.Paragraphs.Item(86)
$opmaak.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item(87)
$p1.Range.ListFormat.ListLevelNumber = 2
$p1.Range.Text = "No kwallen if no stroom goal"

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item(88)
$p2.Range.ListFormat.ListLevelNumber = 2
$p2.Range.Text = "Change win screen (delay)"

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(89)
$p3.Style = "Normal"
try {
    $p3.Range.Delete()
    Write-Host "delete succeeded"
} catch {
    Write-Host "delete FAILED: $_"
}

Write-Host "step1 done"
